$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "she has already taken" -> "they have already taken"
#    "she has not" -> "they have not"
#    "with her in any way" -> "with them in any way"
#    (search text is kept within the existing run boundaries, i.e. it
#     stops right before the separate "exam" runs, so those runs are
#     left untouched). Do this before the "her" -> "their" swap below so
#     that only one whole-word "her" remains in the document.
# ------------------------------------------------------------------
$d.Content.Find.Execute(", I will make sure that she has already taken the ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", I will make sure that they have already taken the ", 2) | Out-Null

$d.Content.Find.Execute(".  If she has not, I will not communicate with her in any way about the ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ".  If they have not, I will not communicate with them in any way about the ", 2) | Out-Null

# ------------------------------------------------------------------
# 2. "her" -> "their"  (within "...other than the Instructor or her assistants.")
#    Match whole word only, so it targets just that remaining "her".
# ------------------------------------------------------------------
$d.Content.Find.Execute("her", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "their", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Paragraph indent: add right indent of -432 twips (-21.6 pt) to the
#    paragraph about internet usage / CSSE 120 web site / GitHub.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*use the internet beyond*") {
        $p.Format.RightIndent = -21.6
        break
    }
}

# ------------------------------------------------------------------
# 4. Insert a manual line break between "the" and "CSSE 120 web site"
# ------------------------------------------------------------------
$d.Content.Find.Execute("directly from the CSSE 120 web site", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "directly from the^lCSSE 120 web site", 2) | Out-Null

# ------------------------------------------------------------------
# 5. ", and my own GitHub repositories." -> " and my own code repository for CSSE 120."
# ------------------------------------------------------------------
$d.Content.Find.Execute(", and my own GitHub repositories.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " and my own code repository for CSSE 120.", 2) | Out-Null

# ------------------------------------------------------------------
# 6. "Academic Honesty" -> "Academic Integrity" (both checkbox paragraphs;
#    Replace:=2 / wdReplaceAll handles every occurrence in one call)
# ------------------------------------------------------------------
$d.Content.Find.Execute("Academic Honesty", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Academic Integrity", 2) | Out-Null

# ------------------------------------------------------------------
# 7. Remove the _GoBack bookmark at the end of the document.
#    (_GoBack is a hidden bookmark, so address it by name directly
#    rather than relying on enumeration.)
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
